# Insert a new weekly record at row 144 for "Vega Monumental Concepción" / "Espinaca".
# Every existing row from 144 down to 151 shifts down by one (to 145-152); this
# script inserts a blank row at 144 (which Excel does automatically, pushing the
# rest down) and then fills the new row 144 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 144 - everything below shifts down.
$ws.Rows.Item(144).Insert()

# Copy the date column's number format from the row just below (now row 145,
# which held the data that used to be in row 144) so the new date cell matches
# the existing "Fecha" column formatting.
$ws.Cells.Item(144, 4).NumberFormat = $ws.Cells.Item(145, 4).NumberFormat

# Populate the new row 144 with the latest weekly price record.
$ws.Cells.Item(144, 1).Value  = 11
$ws.Cells.Item(144, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(144, 3).Value  = "Bíobío"
$ws.Cells.Item(144, 4).Value  = 45267
$ws.Cells.Item(144, 5).Value  = 8
$ws.Cells.Item(144, 6).Value  = 100112012
$ws.Cells.Item(144, 7).Value  = "Espinaca"
$ws.Cells.Item(144, 8).Value  = "Sin especificar"
$ws.Cells.Item(144, 9).Value  = "Primera"
$ws.Cells.Item(144, 10).Value = 50
$ws.Cells.Item(144, 11).Value = 6500
$ws.Cells.Item(144, 12).Value = 6500
$ws.Cells.Item(144, 13).Value = 6500
$ws.Cells.Item(144, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(144, 15).Value = "Región Metropolitana"
$ws.Cells.Item(144, 16).Value = 650
$ws.Cells.Item(144, 17).Value = 10
$ws.Cells.Item(144, 18).Value = "Hortaliza"
